$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the workbook's default/Normal font from Calibri to Arial.
$wb.Styles("Normal").Font.Name = "Arial"

# Header row.
$ws.Range("A1").Value = "Level"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Code"

# First data row: level number, time, and a text "code" value that must
# stay text (it starts with '+' so it needs the apostrophe text-prefix).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2.87
$ws.Range("C2").Value = "'+ACAtCAAgMA5wAkADQCMABxAkADQQMAJwAkEDQSMAhwAkADQQMAJwAkCDQQMAJwAkFDQCMAZx'AEGDQKMApwAEiAoCMApwAkCDQKMApwAkCDQKMAB"

# Remaining level numbers, column A only, rows 3-16 (levels 2-15).
for ($level = 2; $level -le 15; $level++) {
    $row = $level + 1
    $ws.Cells.Item($row, 1).Value = $level
}

[void]$ws.Range("C2").Select()
